# support _null_ as a valid cell entry
#
# Adds a new "someNulls" demonstration row (testing that the parser
# accepts the literal "_null_" token) to both the BasicParsing sheet
# (as a new row 17, pushing the trailing rows down) and to the
# .EmptyRowTesting sheet (as a new row 6).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# BasicParsing: insert a new row 17 with the someNulls sample data
# ---------------------------------------------------------------
$wsBasic = $wb.Worksheets.Item("BasicParsing")

$wsBasic.Rows.Item(17).Insert()

$wsBasic.Range("A17").Value = "someNulls"
$wsBasic.Range("B17").Value = 1
$wsBasic.Range("C17").Value = "_null_"
$wsBasic.Range("D17").Value = $true
$wsBasic.Range("E17").Value = "_null_"
$wsBasic.Range("F17").Value = "hola"
$wsBasic.Range("G17").Value = '"all in"'
$wsBasic.Range("H17").Value = "_null_"
$wsBasic.Rows.Item(17).RowHeight = 18

# widen column H so the new "_null_" values are fully visible
$wsBasic.Columns.Item(8).ColumnWidth = 35.142857142857146

# ---------------------------------------------------------------
# .EmptyRowTesting: mirror the same sample row down at row 6,
# copying BasicParsing's formatting for the new row
# ---------------------------------------------------------------
$wsEmpty = $wb.Worksheets.Item(".EmptyRowTesting")

$wsBasic.Range("A17:H17").Copy()
$wsEmpty.Range("A6").PasteSpecial(-4122)

$wsEmpty.Range("A6").Value = "someNulls"
$wsEmpty.Range("B6").Value = 1
$wsEmpty.Range("C6").Value = "_null_"
$wsEmpty.Range("D6").Value = $true
$wsEmpty.Range("E6").Value = "_null_"
$wsEmpty.Range("F6").Value = "hola"
$wsEmpty.Range("G6").Value = '"all in"'
$wsEmpty.Range("H6").Value = "_null_"
$wsEmpty.Rows.Item(6).RowHeight = 18

[void]$wsEmpty.Range("E16").Select()

# ---------------------------------------------------------------
# Restore BasicParsing as the active sheet/tab, with the new row
# selected (whole-row selection, as left by the editing session).
# This also clears the tabSelected flag that used to sit on
# ErrorCasesParsing, since only one sheet can carry it.
# ---------------------------------------------------------------
$wsBasic.Activate()
[void]$wsBasic.Rows.Item(17).Select()
